$wb = $excel.ActiveWorkbook

# Sheet 1: OUTP1M_RATIO
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = 15.3968702654885
$ws1.Range("A2").Value = 10.7628517865512
$ws1.Range("A3").Value = 1.95950731912582
$ws1.Range("A4").Value = 3.26581135713968
$ws1.Range("A5").Value = 1.18779324855919
$ws1.Range("A6").Value = 1.02538651360578
$ws1.Range("A7").Value = 1.11386708592817
$ws1.Range("A8").Value = 1.12750234811552
$ws1.Range("A9").Value = 1.05605728207215

# Sheet 2: CHRONIC_RATIO
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = 15.5966125456518
$ws2.Range("A2").Value = 11.1898032521172
$ws2.Range("A3").Value = 2.20991141563996
$ws2.Range("A4").Value = 1.10663384212581
$ws2.Range("A5").Value = 1.05294349186372
$ws2.Range("A6").Value = 2.81454112385297
$ws2.Range("A7").Value = 1.04882038640453
$ws2.Range("A8").Value = 1.02175818921065
$ws2.Range("A9").Value = 1.35385937387962
$ws2.Range("A10").Value = 1.2807555387255
